$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 90-133 (1-indexed, matching worksheet row numbers) have their
# "BH.NA" (col G) and "qual.pval" (col H) values replaced with "NA",
# mirroring what already sits in row 134 onward.
$ws.Range("G90:H133").Value = "NA"

# Reflect the author's final selection (G90:H134 highlighted).
$ws.Range("G90:H134").Select()
